$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 8 — this pushes the existing rows
# 8..75 down to 9..76 (matching the diff, where row N's old content
# now lives at row N+1, and the former row-75 data lands on row 76).
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new record.
$ws.Cells.Item(8, 1).Value = 1
$ws.Cells.Item(8, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(8, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(8, 4).Value = 44490
$ws.Cells.Item(8, 5).Value = 15
$ws.Cells.Item(8, 6).Value = "Fruta"
$ws.Cells.Item(8, 7).Value = 100102
$ws.Cells.Item(8, 8).Value = "Cítricos"
$ws.Cells.Item(8, 9).Value = 100102004
$ws.Cells.Item(8, 10).Value = "Mandarina"
$ws.Cells.Item(8, 11).Value = "Murcott"
$ws.Cells.Item(8, 12).Value = "Segunda"
$ws.Cells.Item(8, 13).Value = 250
$ws.Cells.Item(8, 14).Value = 14000
$ws.Cells.Item(8, 15).Value = 15000
$ws.Cells.Item(8, 16).Value = 14500
$ws.Cells.Item(8, 17).Value = "`$/caja 20 kilos"
$ws.Cells.Item(8, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(8, 19).Value = 725
$ws.Cells.Item(8, 20).Value = 20
